$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 3 (EdwinMiranda)
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 100
$ws.Range("G3").Value = 1
$ws.Range("H3").Value = 100
$ws.Range("Q3").Value = 0.71

# Row 4 (johan)
$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 100
$ws.Range("Q4").Value = 0.71

# Row 7 (Nancy Moreno)
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("Q7").Value = 0.71

# Row 9 (NicolasTovar)
$ws.Range("M9").Value = 2
$ws.Range("N9").Value = 66.67
